# Add a new forecast column "BB" (next period, date 45986) to the existing
# YoY component forecast table. Most rows carry forward the last known
# forecast value from column BA, while the most recent rows (19-21) get a
# freshly revised forecast.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy BA1 (date header, with its style/formatting) into BB1, then overwrite
# with the new period's date serial value.
$ws.Range("BA1").Copy($ws.Range("BB1"))
$ws.Range("BB1").Value = 45986

# Rows 3-18: new column simply repeats the last forecast value held in BA.
$carryRows = 3..18
foreach ($r in $carryRows) {
    $srcCell = $ws.Range("BA$r")
    $dstCell = $ws.Range("BB$r")
    $dstCell.Value = $srcCell.Value()
}

# Rows 19-21: updated/revised forecast values for the newest data points.
$ws.Range("BB19").Value = 2.622852459381209
$ws.Range("BB20").Value = 2.689750575689809
$ws.Range("BB21").Value = 3.108260574154809
